# Client form edit:
# - reorder the header/data columns of the "Clientes" sheet
# - drop the usuario_cipre / contrasena columns (no longer exposed in the editor)
# - keep every existing header + its row-2 sample value glued together while moving
#
# NOTE: we intentionally move cells with Range.Cut(destination) instead of
# re-typing literal values into new cells. Cut preserves both the original
# cell style (so header cells keep their bold/bordered look) and the original
# stored type (so numeric-looking text such as "5551234567" / "2025-12-11"
# stays text instead of being re-interpreted as a number/date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old column letter -> new column letter, in the new left-to-right order.
$moves = @(
    @{Old="A";  New="A"}   # id
    @{Old="P";  New="B"}   # tipo_tramite
    @{Old="AB"; New="C"}   # producto
    @{Old="O";  New="D"}   # fuente
    @{Old="AC"; New="E"}   # fuente_base
    @{Old="B";  New="F"}   # nombre
    @{Old="L";  New="G"}   # telefono
    @{Old="C";  New="H"}   # sucursal
    @{Old="Q";  New="I"}   # capacidad
    @{Old="I";  New="J"}   # monto_final
    @{Old="H";  New="K"}   # monto_propuesta
    @{Old="R";  New="L"}   # plazo
    @{Old="S";  New="M"}   # estado_civil
    @{Old="T";  New="N"}   # tipo_vivienda
    @{Old="M";  New="O"}   # correo
    @{Old="U";  New="P"}   # ref1_nombre
    @{Old="V";  New="Q"}   # ref1_telefono
    @{Old="W";  New="R"}   # ref1_parentesco
    @{Old="X";  New="S"}   # ref2_nombre
    @{Old="Y";  New="T"}   # ref2_telefono
    @{Old="Z";  New="U"}   # ref2_parentesco
    @{Old="AA"; New="V"}   # antiguedad_cuenta
    @{Old="D";  New="W"}   # asesor
    @{Old="E";  New="X"}   # fecha_ingreso
    @{Old="F";  New="Y"}   # fecha_dispersion
    @{Old="G";  New="Z"}   # estatus
    @{Old="J";  New="AA"}  # observaciones
    @{Old="K";  New="AB"}  # score
    @{Old="N";  New="AC"}  # analista
)

# ---------------------------------------------------------------------------
# Phase 1: move every (header, value) column pair out to a staging area
# (columns far to the right) using the NEW column order, so that later moves
# never clobber a source column that hasn't been relocated yet.
#
# Row 2 is sparse in the original sheet (most of the new/blank columns have
# no sample value at all), so we only Cut the data cell when it actually
# holds something - Cut-ing an empty cell would otherwise materialise a new
# (empty) <c> element at the destination and leave the sheet's used range
# bigger than it should be.
# ---------------------------------------------------------------------------
$stageStartCol = 200   # arbitrary far-away staging column (column "GR")
for ($i = 0; $i -lt $moves.Count; $i++) {
    $oldCol = $moves[$i].Old
    $stageCol = $stageStartCol + $i

    $headerSrc = $ws.Range("$($oldCol)1")
    $headerDst = $ws.Cells.Item(1, $stageCol)
    $headerSrc.Cut($headerDst)

    $valueSrc = $ws.Range("$($oldCol)2")
    if ($valueSrc.Value2 -ne $null -and $valueSrc.Value2 -ne "") {
        $valueDst = $ws.Cells.Item(2, $stageCol)
        $valueSrc.Cut($valueDst)
    }
}

# ---------------------------------------------------------------------------
# Phase 2: clear whatever remains in the old layout (this drops the
# usuario_cipre / contrasena columns for good) and reset the used range.
# ---------------------------------------------------------------------------
$ws.Range("A1:AE2").Clear()

# ---------------------------------------------------------------------------
# Phase 3: move everything from staging back into its final position.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt $moves.Count; $i++) {
    $newCol = $moves[$i].New
    $stageCol = $stageStartCol + $i

    $headerSrc = $ws.Cells.Item(1, $stageCol)
    $headerDst = $ws.Range("$($newCol)1")
    $headerSrc.Cut($headerDst)

    $valueSrc = $ws.Cells.Item(2, $stageCol)
    if ($valueSrc.Value2 -ne $null -and $valueSrc.Value2 -ne "") {
        $valueDst = $ws.Range("$($newCol)2")
        $valueSrc.Cut($valueDst)
    }
}

# ---------------------------------------------------------------------------
# Phase 4: scrub any formatting remnants left behind in the staging area so
# the sheet's dimension/used-range shrinks back down to the real data.
# ---------------------------------------------------------------------------
$lastStageCol = $stageStartCol + $moves.Count - 1
$stageAll = $ws.Range($ws.Cells.Item(1, $stageStartCol), $ws.Cells.Item(2, $lastStageCol))
$stageAll.Clear()
